$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 3 through 6, leaving only rows 1-2
$ws.Range("A3:C6").EntireRow.Delete() | Out-Null

# Update the remaining values to match the new data
$ws.Range("A1").Value = 166
$ws.Range("B1").Value = 166
$ws.Range("C1").Value = 198

$ws.Range("A2").Value = 166
$ws.Range("B2").Value = 210.3999999999996
$ws.Range("C2").Value = 166
